# Adds a new "Busy Street" sound-reference entry (bulleted list item with a
# hyperlink to the freesound.org page) right after the existing "BBC: ..."
# list item, matching the "Added sounds array + busyStreet.wav" commit.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "BBC: ..." bullet - it is the last populated list item.
# ------------------------------------------------------------------
$bbcIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "BBC:*") {
        $bbcIndex = $i
    }
}

$bbcPara = $d.Paragraphs.Item($bbcIndex)

# The old "_GoBack" bookmark Word leaves around the last edit position is
# stale after this edit - drop it (matches the target document, which no
# longer carries it).
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
}

# ------------------------------------------------------------------
# Insert a brand new list paragraph right after the BBC bullet. The new
# paragraph mark inherits BBC's pPr (pStyle "ListParagraph" + numPr
# ilvl=0/numId=1), so the bullet list formatting carries over for free.
# ------------------------------------------------------------------
$bbcPara.Range.InsertParagraphAfter()
$newIndex = $bbcIndex + 1
$newPara = $d.Paragraphs.Item($newIndex)

# Label text.
$labelRange = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$labelRange.InsertAfter("Busy Street: ")

# URL text + hyperlink.
$url = "https://freesound.org/people/Rico_Casazza/sounds/538952/"
$urlRange = $d.Range($labelRange.End, $labelRange.End)
$urlRange.InsertAfter($url)
$d.Hyperlinks.Add($urlRange, $url)

# Trailing space, also wrapped in its own hyperlink (mirrors the pattern
# already used for the other entries, e.g. the Seagull / Dog bark bullets,
# where the space after the link text got swept into the hyperlink too).
$afterLinkPos = $d.Paragraphs.Item($newIndex).Range.End - 1
$spaceRange = $d.Range($afterLinkPos, $afterLinkPos)
$spaceRange.InsertAfter(" ")

$newEnd = $d.Paragraphs.Item($newIndex).Range.End
$spaceCharRange = $d.Range($newEnd - 2, $newEnd - 1)
$d.Hyperlinks.Add($spaceCharRange, $url)
